$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to Text format so numeric-looking values
# (e.g. "311.27") are not auto-converted into floating point numbers,
# matching the original inline-string cell semantics.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "41.846.17"
$ws.Range("E2").Value = "  -1.87%  "
# Row 3
$ws.Range("D3").Value = "2.301.54"
$ws.Range("E3").Value = "  -2.39%  "
# Row 4
$ws.Range("E4").Value = "  -0.22%  "
# Row 5
$ws.Range("D5").Value = "311.27"
$ws.Range("E5").Value = "  -6.28%  "
# Row 6
$ws.Range("D6").Value = "104.99"
$ws.Range("E6").Value = "  +4.66%  "
# Row 7
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -1.93%  "
# Row 8
$ws.Range("E8").Value = "  -0.08%  "
# Row 9
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -4.31%  "
# Row 10
$ws.Range("D10").Value = "39.96"
$ws.Range("E10").Value = "  +0.41%  "
# Row 11
$ws.Range("D11").Value = "0.0911"
$ws.Range("E11").Value = "  -1.25%  "
# Row 12
$ws.Range("D12").Value = "8.27"
$ws.Range("E12").Value = "  -1.99%  "
# Row 13
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  -0.03%  "
# Row 14
$ws.Range("D14").Value = "0.972"
$ws.Range("E14").Value = "  -2.33%  "
# Row 15
$ws.Range("D15").Value = "15.41"
$ws.Range("E15").Value = "  -5.28%  "
# Row 16
$ws.Range("D16").Value = "2.642.03"
$ws.Range("E16").Value = "  -2.80%  "
# Row 17
$ws.Range("D17").Value = "2.295.41"
$ws.Range("E17").Value = "  -2.51%  "
# Row 18
$ws.Range("D18").Value = "41.868.03"
$ws.Range("E18").Value = "  -1.76%  "
# Row 19
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  -1.33%  "
# Row 20
$ws.Range("E20").Value = "  -2.36%  "
# Row 21
$ws.Range("D21").Value = "74.37"
$ws.Range("E21").Value = "  -1.85%  "
# Row 22
$ws.Range("D22").Value = "3.49"
$ws.Range("E22").Value = "  -7.46%  "
# Row 23
$ws.Range("D23").Value = "258.31"
$ws.Range("E23").Value = "  -4.36%  "
# Row 24
$ws.Range("E24").Value = "  -2.63%  "
# Row 25
$ws.Range("D25").Value = "9.23"
$ws.Range("E25").Value = "  -7.51%  "
# Row 26
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.66%  "
# Row 27
$ws.Range("D27").Value = "10.97"
$ws.Range("E27").Value = "  -4.19%  "
# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "22.80"
$ws.Range("E28").Value = "  -1.87%  "
# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +0.70%  "
# Row 30
$ws.Range("D30").Value = "35.40"
$ws.Range("E30").Value = "  -0.31%  "
# Row 31
$ws.Range("D31").Value = "162.40"
$ws.Range("E31").Value = "  -7.33%  "
# Row 32
$ws.Range("D32").Value = "0.0885"
$ws.Range("E32").Value = "  -2.07%  "
# Row 33
$ws.Range("D33").Value = "2.90"
$ws.Range("E33").Value = "  -6.42%  "
# Row 34
$ws.Range("D34").Value = "5.81"
$ws.Range("E34").Value = "  -4.17%  "
# Row 35
$ws.Range("D35").Value = "0.129"
$ws.Range("E35").Value = "  -3.36%  "
# Row 36
$ws.Range("E36").Value = "  +9.15%  "
# Row 37
$ws.Range("D37").Value = "4.50"
$ws.Range("E37").Value = "  -1.98%  "
# Row 38
$ws.Range("D38").Value = "0.0350"
$ws.Range("E38").Value = "  -2.46%  "
# Row 39
$ws.Range("D39").Value = "3.64"
$ws.Range("E39").Value = "  -4.45%  "
# Row 40
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  -6.87%  "
# Row 41
$ws.Range("D41").Value = "99.12"
$ws.Range("E41").Value = "  +9.39%  "
# Row 42
$ws.Range("D42").Value = "1.47"
$ws.Range("E42").Value = "  -3.66%  "
# Row 43
$ws.Range("D43").Value = "69.61"
$ws.Range("E43").Value = "  -0.69%  "
# Row 44
$ws.Range("D44").Value = "0.229"
$ws.Range("E44").Value = "  -2.29%  "
# Row 45
$ws.Range("E45").Value = "  -0.29%  "
# Row 46
$ws.Range("D46").Value = "12.06"
$ws.Range("E46").Value = "  +0.64%  "
# Row 47
$ws.Range("D47").Value = "111.53"
$ws.Range("E47").Value = "  -5.35%  "
# Row 48
$ws.Range("D48").Value = "5.35"
$ws.Range("E48").Value = "  -2.15%  "
# Row 49
$ws.Range("D49").Value = "8.90"
$ws.Range("E49").Value = "  -2.02%  "
# Row 50
$ws.Range("D50").Value = "73.73"
$ws.Range("E50").Value = "  +5.59%  "
# Row 51
$ws.Range("E51").Value = "  -0.49%  "
